# The document has several runs whose visible text is exactly " n"
# (a space followed by the letter n). Each of them sits right after a
# run with identical character formatting (same w:rPr), so a naive
# Range.Text / Find-Replace assignment causes the host to re-merge the
# freshly edited run back together with its neighbouring, identically
# formatted runs (losing the run boundaries the target XML wants to
# keep). To stop that happening we briefly give the target run a
# distinguishing direct-formatting toggle (Bold) before rewriting its
# text, then switch that toggle back off once the new text is in
# place - by that point the run's formatting is identical to the
# original again, but the automatic "merge identically formatted runs"
# pass has already happened against the *old* (distinct) formatting, so
# the run is left standing on its own with the new text, just like the
# target XML expects.

$d = $word.ActiveDocument

$old = " n"
$new = "tjjt"

$searchStart = 0
for ($i = 0; $i -lt 50; $i++) {
    $probe = $d.Range($searchStart, $d.Content.End)
    $found = $probe.Find.Execute($old, $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    $matchStart = $probe.Start
    $matchEnd = $probe.End

    $target = $d.Range($matchStart, $matchEnd)
    $target.Font.Bold = $true
    $target.Text = $new

    $newEnd = $matchStart + $new.Length
    $resetRange = $d.Range($matchStart, $newEnd)
    $resetRange.Font.Bold = $false

    $searchStart = $newEnd
}
